$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 651
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 651
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 651
$ws.Cells.Item(12, 13).Value = ""
$ws.Cells.Item(12, 14).Value = -991
$ws.Cells.Item(17, 8).Value = 1581.3704
$ws.Cells.Item(17, 10).Value = 1581.3704
$ws.Cells.Item(17, 12).Value = 4744.1112
$ws.Cells.Item(17, 14).Value = -5080.1112
$ws.Cells.Item(55, 8).Value = 331.25
$ws.Cells.Item(55, 9).Value = 212.6
$ws.Cells.Item(55, 10).Value = 529
$ws.Cells.Item(55, 11).Value = 212.6
$ws.Cells.Item(55, 12).Value = 529
$ws.Cells.Item(55, 13).Value = 1.400000000000006
$ws.Cells.Item(55, 14).Value = -957
$ws.Cells.Item(64, 8).Value = 4807.615
$ws.Cells.Item(64, 9).Value = 4799.8
$ws.Cells.Item(64, 11).Value = 4799.8
$ws.Cells.Item(64, 13).Value = -4551.8
$ws.Cells.Item(67, 8).Value = 4807.615
$ws.Cells.Item(67, 9).Value = 4799.8
$ws.Cells.Item(67, 11).Value = 4799.8
$ws.Cells.Item(67, 13).Value = -3941.8
$ws.Cells.Item(98, 8).Value = 1120.25
$ws.Cells.Item(98, 9).Value = 855.8333
$ws.Cells.Item(98, 11).Value = 855.8333
$ws.Cells.Item(98, 13).Value = 642.1667
$ws.Cells.Item(100, 8).Value = 2090.8333
$ws.Cells.Item(100, 9).Value = 2019
$ws.Cells.Item(100, 11).Value = 2019
$ws.Cells.Item(100, 13).Value = -1478
$ws.Cells.Item(116, 8).Value = 7612.5
$ws.Cells.Item(116, 9).Value = 10000
$ws.Cells.Item(116, 10).Value = 6816.6665
$ws.Cells.Item(116, 11).Value = 10000
$ws.Cells.Item(116, 12).Value = 6816.6665
$ws.Cells.Item(116, 13).Value = -6558
$ws.Cells.Item(116, 14).Value = -13700.6665
$ws.Cells.Item(118, 8).Value = 199.66667
$ws.Cells.Item(118, 9).Value = 199.66667
$ws.Cells.Item(118, 11).Value = 599.00001
$ws.Cells.Item(118, 13).Value = 1057.99999
$ws.Cells.Item(122, 8).Value = 1120.25
$ws.Cells.Item(122, 9).Value = 855.8333
$ws.Cells.Item(122, 11).Value = 2567.4999
$ws.Cells.Item(122, 13).Value = -117.4998999999998
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1431.381
$ws.Cells.Item(32, 9).Value = 1257.4237
$ws.Cells.Item(32, 11).Value = 1257.4237
$ws.Cells.Item(32, 13).Value = -970.4237000000001
$ws.Cells.Item(92, 8).Value = 27249.25
$ws.Cells.Item(92, 10).Value = 27249.25
$ws.Cells.Item(92, 12).Value = 27249.25
$ws.Cells.Item(92, 14).Value = -32241.25
$ws.Cells.Item(122, 8).Value = 4337.3335
$ws.Cells.Item(122, 9).Value = 4006
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 12018
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -9568
$ws.Cells.Item(122, 14).Value = -19900
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 170.33333
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 12).Value = 200
$ws.Cells.Item(22, 14).Value = -546
$ws.Cells.Item(100, 8).Value = 12327.2
$ws.Cells.Item(100, 10).Value = 12327.2
$ws.Cells.Item(100, 12).Value = 12327.2
$ws.Cells.Item(100, 14).Value = -14491.2
$ws.Cells.Item(104, 8).Value = 50684
$ws.Cells.Item(104, 10).Value = 50684
$ws.Cells.Item(104, 12).Value = 50684
$ws.Cells.Item(104, 14).Value = -57672
$ws.Cells.Item(105, 8).Value = 3773.375
$ws.Cells.Item(105, 9).Value = 3964.6667
$ws.Cells.Item(105, 11).Value = 3964.6667
$ws.Cells.Item(105, 13).Value = -2217.6667
$ws.Cells.Item(132, 8).Value = 1820.5385
$ws.Cells.Item(132, 9).Value = 1806
$ws.Cells.Item(132, 11).Value = 5418
$ws.Cells.Item(132, 13).Value = -2888
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 400
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 400
$ws.Cells.Item(22, 13).Value = ""
$ws.Cells.Item(22, 14).Value = -1100
$ws.Cells.Item(58, 8).Value = 3119.4707
$ws.Cells.Item(58, 9).Value = 2633.1
$ws.Cells.Item(58, 10).Value = 3814.2856
$ws.Cells.Item(58, 11).Value = 2633.1
$ws.Cells.Item(58, 12).Value = 3814.2856
$ws.Cells.Item(58, 13).Value = -2430.1
$ws.Cells.Item(58, 14).Value = -4220.2856
$ws.Cells.Item(86, 8).Value = 11715.125
$ws.Cells.Item(86, 10).Value = 10878
$ws.Cells.Item(86, 12).Value = 10878
$ws.Cells.Item(86, 14).Value = -13124
$ws.Cells.Item(89, 8).Value = 11715.125
$ws.Cells.Item(89, 10).Value = 10878
$ws.Cells.Item(89, 12).Value = 54390
$ws.Cells.Item(89, 14).Value = -65622
$ws.Cells.Item(99, 8).Value = 4954
$ws.Cells.Item(99, 9).Value = 3707.6667
$ws.Cells.Item(99, 10).Value = 5888.75
$ws.Cells.Item(99, 11).Value = 3707.6667
$ws.Cells.Item(99, 12).Value = 5888.75
$ws.Cells.Item(99, 13).Value = -2209.6667
$ws.Cells.Item(99, 14).Value = -8884.75
$ws.Cells.Item(102, 8).Value = 45120.5
$ws.Cells.Item(102, 10).Value = 45120.5
$ws.Cells.Item(102, 12).Value = 45120.5
$ws.Cells.Item(102, 14).Value = -49988.5
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = ""
$ws.Cells.Item(105, 14).Value = ""
$ws.Cells.Item(109, 8).Value = 66794.664
$ws.Cells.Item(109, 10).Value = 66794.664
$ws.Cells.Item(109, 12).Value = 66794.664
$ws.Cells.Item(109, 14).Value = -68874.664
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).Value = ""
$ws.Cells.Item(126, 8).Value = 4954
$ws.Cells.Item(126, 9).Value = 3707.6667
$ws.Cells.Item(126, 10).Value = 5888.75
$ws.Cells.Item(126, 11).Value = 11123.0001
$ws.Cells.Item(126, 12).Value = 17666.25
$ws.Cells.Item(126, 13).Value = -8653.000100000001
$ws.Cells.Item(126, 14).Value = -22606.25
$ws.Cells.Item(134, 8).Value = 7806.2
$ws.Cells.Item(134, 9).Value = 7506.9443
$ws.Cells.Item(134, 10).Value = 10499.5
$ws.Cells.Item(134, 11).Value = 22520.8329
$ws.Cells.Item(134, 12).Value = 31498.5
$ws.Cells.Item(134, 13).Value = -19985.8329
$ws.Cells.Item(134, 14).Value = -36568.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 467.85715
$ws.Cells.Item(9, 9).Value = 443.75
$ws.Cells.Item(9, 10).Value = 500
$ws.Cells.Item(9, 11).Value = 1331.25
$ws.Cells.Item(9, 12).Value = 1500
$ws.Cells.Item(9, 13).Value = -1107.25
$ws.Cells.Item(9, 14).Value = -1948
$ws.Cells.Item(92, 8).Value = 278
$ws.Cells.Item(92, 9).Value = 408.5
$ws.Cells.Item(92, 10).Value = 147.5
$ws.Cells.Item(92, 11).Value = 1225.5
$ws.Cells.Item(92, 12).Value = 442.5
$ws.Cells.Item(92, 13).Value = 22.5
$ws.Cells.Item(92, 14).Value = -2938.5
$ws.Cells.Item(119, 8).Value = 20
$ws.Cells.Item(119, 9).Value = 20
$ws.Cells.Item(119, 11).Value = 60
$ws.Cells.Item(119, 13).Value = 4778
$ws.Cells.Item(120, 8).Value = 30
$ws.Cells.Item(120, 9).Value = 30
$ws.Cells.Item(120, 11).Value = 90
$ws.Cells.Item(120, 13).Value = 4748
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).Value = ""
$ws.Cells.Item(134, 8).Value = 2782.6667
$ws.Cells.Item(134, 9).Value = 2868
$ws.Cells.Item(134, 11).Value = 8604
$ws.Cells.Item(134, 13).Value = -6069
$ws.Cells.Item(136, 8).Value = 3119.4707
$ws.Cells.Item(136, 9).Value = 2633.1
$ws.Cells.Item(136, 10).Value = 3814.2856
$ws.Cells.Item(136, 11).Value = 7899.299999999999
$ws.Cells.Item(136, 12).Value = 11442.8568
$ws.Cells.Item(136, 13).Value = -5349.299999999999
$ws.Cells.Item(136, 14).Value = -16542.8568
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3000
$ws.Cells.Item(80, 9).Value = 3000
$ws.Cells.Item(80, 11).Value = 3000
$ws.Cells.Item(80, 13).Value = -2002
$ws.Cells.Item(83, 8).Value = 3000
$ws.Cells.Item(83, 9).Value = 3000
$ws.Cells.Item(83, 11).Value = 15000
$ws.Cells.Item(83, 13).Value = -10008
$ws.Cells.Item(122, 8).Value = 2611.8823
$ws.Cells.Item(122, 9).Value = 2446.6924
$ws.Cells.Item(122, 11).Value = 7340.0772
$ws.Cells.Item(122, 13).Value = -4890.0772
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 499.6
$ws.Cells.Item(16, 9).Value = 499.6
$ws.Cells.Item(16, 11).Value = 499.6
$ws.Cells.Item(16, 13).Value = -329.6
$ws.Cells.Item(39, 8).Value = 25888
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 13).Value = ""
$ws.Cells.Item(46, 8).Value = 3789.4211
$ws.Cells.Item(46, 9).Value = 1750
$ws.Cells.Item(46, 10).Value = 4029.353
$ws.Cells.Item(46, 11).Value = 1750
$ws.Cells.Item(46, 12).Value = 4029.353
$ws.Cells.Item(46, 13).Value = -1562
$ws.Cells.Item(46, 14).Value = -4405.353
$ws.Cells.Item(55, 8).Value = 208.88889
$ws.Cells.Item(55, 9).Value = 209.4
$ws.Cells.Item(55, 10).Value = 208.6923
$ws.Cells.Item(55, 11).Value = 209.4
$ws.Cells.Item(55, 12).Value = 208.6923
$ws.Cells.Item(55, 13).Value = -36.40000000000001
$ws.Cells.Item(55, 14).Value = -554.6922999999999
$ws.Cells.Item(68, 8).Value = 36066.168
$ws.Cells.Item(68, 9).Value = 1400
$ws.Cells.Item(68, 10).Value = 42999.4
$ws.Cells.Item(68, 11).Value = 1400
$ws.Cells.Item(68, 12).Value = 42999.4
$ws.Cells.Item(68, 13).Value = -651
$ws.Cells.Item(68, 14).Value = -44497.4
$ws.Cells.Item(71, 8).Value = 36066.168
$ws.Cells.Item(71, 9).Value = 1400
$ws.Cells.Item(71, 10).Value = 42999.4
$ws.Cells.Item(71, 11).Value = 7000
$ws.Cells.Item(71, 12).Value = 214997
$ws.Cells.Item(71, 13).Value = -3256
$ws.Cells.Item(71, 14).Value = -222485
$ws.Cells.Item(93, 8).Value = 833.1667
$ws.Cells.Item(93, 9).Value = 783.3333
$ws.Cells.Item(93, 10).Value = 883
$ws.Cells.Item(93, 11).Value = 783.3333
$ws.Cells.Item(93, 12).Value = 883
$ws.Cells.Item(93, 13).Value = 464.6667
$ws.Cells.Item(93, 14).Value = -3379
$ws.Cells.Item(132, 8).Value = 2531.8462
$ws.Cells.Item(132, 9).Value = 2601.6365
$ws.Cells.Item(132, 10).Value = 2148
$ws.Cells.Item(132, 11).Value = 7804.9095
$ws.Cells.Item(132, 12).Value = 6444
$ws.Cells.Item(132, 13).Value = -5274.9095
$ws.Cells.Item(132, 14).Value = -11504
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 1891.6666
$ws.Cells.Item(4, 9).Value = 526
$ws.Cells.Item(4, 10).Value = 2574.5
$ws.Cells.Item(4, 11).Value = 526
$ws.Cells.Item(4, 12).Value = 2574.5
$ws.Cells.Item(4, 13).Value = -413
$ws.Cells.Item(4, 14).Value = -2800.5
$ws.Cells.Item(54, 8).Value = 22517.5
$ws.Cells.Item(136, 8).Value = 2977.25
$ws.Cells.Item(136, 9).Value = 2603.3333
$ws.Cells.Item(136, 11).Value = 7809.999899999999
$ws.Cells.Item(136, 13).Value = -5259.999899999999
